$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '55.057.57'
$ws.Range('E2').Value = '  +2.37%  '
$ws.Range('D3').Value = '2.277.21'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '''506.79'
$ws.Range('E5').Value = '  +2.84%  '
$ws.Range('D6').Value = '''128.60'
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('D8').Value = '''0.529'
$ws.Range('E8').Value = '  +1.03%  '
$ws.Range('D9').Value = '2.291.48'
$ws.Range('E9').Value = '  +1.83%  '
$ws.Range('D10').Value = '''0.0987'
$ws.Range('E10').Value = '  +4.28%  '
$ws.Range('D11').Value = '''0.155'
$ws.Range('E11').Value = '  +1.07%  '
$ws.Range('E12').Value = '  +6.95%  '
$ws.Range('E13').Value = '  +2.74%  '
$ws.Range('D14').Value = '''23.60'
$ws.Range('E14').Value = '  +4.87%  '
$ws.Range('D15').Value = '2.683.44'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').Value = '55.097.28'
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('E17').Value = '  +1.77%  '
$ws.Range('D18').Value = '2.279.42'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('D19').Value = '''10.38'
$ws.Range('E19').Value = '  +1.76%  '
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('D21').Value = '''314.22'
$ws.Range('E21').Value = '  +4.25%  '
$ws.Range('D22').Value = '''6.56'
$ws.Range('E22').Value = '  +4.58%  '
$ws.Range('D23').Value = '''0.997'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '''59.88'
$ws.Range('E24').Value = '  -1.53%  '
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('E26').Value = '  +5.07%  '
$ws.Range('E27').Value = '  +4.36%  '
$ws.Range('D28').Value = '''171.20'
$ws.Range('E28').Value = '  -0.32%  '
$ws.Range('E29').Value = '  +4.62%  '
$ws.Range('E30').Value = '  +3.21%  '
$ws.Range('E31').Value = '  +3.30%  '
$ws.Range('E32').Value = '  +7.91%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').Value = '''18.00'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('E35').Value = '  -0.50%  '
$ws.Range('D36').Value = '''1.24'
$ws.Range('E36').Value = '  +4.27%  '
$ws.Range('D37').Value = '''0.902'
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('E38').Value = '  +5.73%  '
$ws.Range('D39').Value = '''36.87'
$ws.Range('E39').Value = '  +3.06%  '
$ws.Range('D40').Value = '''1.45'
$ws.Range('E40').Value = '  +4.97%  '
$ws.Range('D41').Value = '''0.375'
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('D42').Value = '''136.59'
$ws.Range('E42').Value = '  +9.82%  '
$ws.Range('D43').Value = '''3.48'
$ws.Range('E43').Value = '  +4.14%  '
$ws.Range('E44').Value = '  +2.48%  '
$ws.Range('D45').Value = '''259.12'
$ws.Range('E45').Value = '  +9.06%  '
$ws.Range('D46').Value = '''0.0507'
$ws.Range('E46').Value = '  +3.92%  '
$ws.Range('E47').Value = '  +3.72%  '
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('E49').Value = '  +4.82%  '
$ws.Range('E50').Value = '  +1.38%  '
$ws.Range('D51').Value = '''16.53'
$ws.Range('E51').Value = '  +3.42%  '
